$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 347, shifting rows 347:373 down to 348:374
$ws.Rows("347:347").Insert()

# Populate the newly inserted row 347 with the new weekly data point
$ws.Cells.Item(347, 1).Value = 9
$ws.Cells.Item(347, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(347, 3).Value = "Metropolitana"
$ws.Cells.Item(347, 4).Value = 44585
$ws.Cells.Item(347, 5).Value = 13
$ws.Cells.Item(347, 6).Value = 100112031
$ws.Cells.Item(347, 7).Value = "Poroto verde"
$ws.Cells.Item(347, 8).Value = "Magnum"
$ws.Cells.Item(347, 9).Value = "Primera"
$ws.Cells.Item(347, 10).Value = 43
$ws.Cells.Item(347, 11).Value = 30000
$ws.Cells.Item(347, 12).Value = 32000
$ws.Cells.Item(347, 13).Value = 31023
$ws.Cells.Item(347, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(347, 15).Value = "Región Metropolitana"
$ws.Cells.Item(347, 16).Value = 1241
$ws.Cells.Item(347, 17).Value = 25
$ws.Cells.Item(347, 18).Value = "Hortaliza"

# Match the date style used by the rest of column D (style index 2, numFmt 165)
$ws.Cells.Item(347, 4).NumberFormat = $ws.Cells.Item(348, 4).NumberFormat
